$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-05-13 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-14 Wednesday", 2) | Out-Null

# Update each math-problem cell in the table, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "72+16=",
    "63-2=",
    "22+7=",
    "6+36=",
    "46-1=",
    "51+16=",
    "19+24=",
    "64-26=",
    "3+40=",
    "49+25=",
    "12+27=",
    "63-40=",
    "76-1=",
    "74-47=",
    "94-60=",
    "42+46=",
    "13-0=",
    "81-45=",
    "55-48=",
    "61-19=",
    "97-38=",
    "72-0=",
    "20+42=",
    "57+12=",
    "72-69=",
    "63-13=",
    "50-2=",
    "39+29=",
    "6+67=",
    "96-73=",
    "33-6=",
    "34+10=",
    "19-5=",
    "38+16=",
    "80-37=",
    "80+6=",
    "65-26=",
    "97-6=",
    "54-52=",
    "28+55=",
    "81-80=",
    "42+37=",
    "59+38=",
    "75-39=",
    "78-2=",
    "22+67=",
    "93-89=",
    "5+32=",
    "58-40=",
    "7+9=",
    "26+41=",
    "60-8=",
    "67-8=",
    "91-57=",
    "63-40=",
    "60+8=",
    "2+56=",
    "31+27=",
    "38+4=",
    "92-54=",
    "81-11=",
    "49+15=",
    "66+30=",
    "14+50=",
    "34-1=",
    "66-50=",
    "28+71=",
    "35-14=",
    "98-63=",
    "60-52=",
    "48-41=",
    "99-10=",
    "3+70=",
    "17+12=",
    "82+4=",
    "25+13=",
    "63+33=",
    "44+54=",
    "67+25=",
    "96-23=",
    "87-62=",
    "15+78=",
    "2+24=",
    "82-74=",
    "99-34=",
    "83-38=",
    "93+6=",
    "42+19=",
    "60+15=",
    "33-13=",
    "96-72=",
    "37+15=",
    "95-80=",
    "78-75=",
    "80+6=",
    "14+4=",
    "10+69=",
    "22-3=",
    "4+73=",
    "87-8="
)

$cols = $t.Columns.Count
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}
